$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.300649642944336
$ws.Range("B1").Value = 4.524884700775146
$ws.Range("C1").Value = 2.646209239959717
$ws.Range("D1").Value = 2.366683959960938
$ws.Range("E1").Value = 2.232384920120239
